$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad) from 45172 to 45175 for existing rows 2:303
$ws.Range("C2:C303").Value = 45175

# 2. Row 303 picks up an explicit row height (matches the rest of the sheet)
$ws.Rows.Item(303).RowHeight = 15

# 3. Append new row 304
$ws.Cells.Item(304, 1).Value = "A 40989-2023"
$ws.Cells.Item(304, 2).Value = 45173
$ws.Cells.Item(304, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(304, 3).Value = 45175
$ws.Cells.Item(304, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(304, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(304, 5).Value = "SÖDERHAMN"
$ws.Cells.Item(304, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(304, 7).Value = 0.9
$ws.Cells.Item(304, 8).Value = 0
$ws.Cells.Item(304, 9).Value = 0
$ws.Cells.Item(304, 10).Value = 0
$ws.Cells.Item(304, 11).Value = 0
$ws.Cells.Item(304, 12).Value = 0
$ws.Cells.Item(304, 13).Value = 0
$ws.Cells.Item(304, 14).Value = 0
$ws.Cells.Item(304, 15).Value = 0
$ws.Cells.Item(304, 16).Value = 0
$ws.Cells.Item(304, 17).Value = 0
$ws.Cells.Item(304, 18).WrapText = $true
$ws.Rows.Item(304).RowHeight = 15

# 4. Append new row 305
$ws.Cells.Item(305, 1).Value = "A 40987-2023"
$ws.Cells.Item(305, 2).Value = 45173
$ws.Cells.Item(305, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(305, 3).Value = 45175
$ws.Cells.Item(305, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(305, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(305, 5).Value = "SÖDERHAMN"
$ws.Cells.Item(305, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(305, 7).Value = 1.7
$ws.Cells.Item(305, 8).Value = 0
$ws.Cells.Item(305, 9).Value = 0
$ws.Cells.Item(305, 10).Value = 0
$ws.Cells.Item(305, 11).Value = 0
$ws.Cells.Item(305, 12).Value = 0
$ws.Cells.Item(305, 13).Value = 0
$ws.Cells.Item(305, 14).Value = 0
$ws.Cells.Item(305, 15).Value = 0
$ws.Cells.Item(305, 16).Value = 0
$ws.Cells.Item(305, 17).Value = 0
$ws.Cells.Item(305, 18).WrapText = $true
